# Update cryptos.xlsx with the latest crypto price/volume snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # The "Price" column stores values that look numeric (e.g. "1.00") but
    # must stay literal text, like the rest of the sheet (inline strings).
    # Writing the digits directly would make Excel auto-convert the cell to
    # a real number, so prefix with an apostrophe (same as a user typing
    # '1.00 into a cell) to force text, then clear the resulting quote-prefix
    # style so no visible formatting changes stick around.
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Rows 40-41: Stellar and Monero swapped ranking position
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D40") "122.70"
$ws.Range("E40").Value = "  +5.16%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D41") "0.112"
$ws.Range("E41").Value = "  +0.56%  "

# Remaining price/volume updates
$ws.Range("D2").Value = '48.177.02'
$ws.Range("E2").Value = '  +1.86%  '
$ws.Range("D3").Value = '2.509.93'
$ws.Range("E3").Value = '  +0.85%  '
Set-TextValue $ws.Range("D4") '1.00'
$ws.Range("E4").Value = '  -0.01%  '
Set-TextValue $ws.Range("D5") '320.68'
$ws.Range("E5").Value = '  -0.07%  '
Set-TextValue $ws.Range("D6") '108.78'
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("E7").Value = '  +1.29%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +1.17%  '
$ws.Range("E10").Value = '  +1.14%  '
Set-TextValue $ws.Range("D11") '20.27'
$ws.Range("E11").Value = '  +10.59%  '
$ws.Range("E12").Value = '  +0.95%  '
Set-TextValue $ws.Range("D13") '0.124'
$ws.Range("E13").Value = '  +0.71%  '
$ws.Range("E14").Value = '  +0.81%  '
$ws.Range("D15").Value = '2.901.73'
$ws.Range("E15").Value = '  +0.70%  '
$ws.Range("D16").Value = '2.503.60'
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("E17").Value = '  +0.23%  '
$ws.Range("D18").Value = '48.024.34'
$ws.Range("E18").Value = '  +1.72%  '
Set-TextValue $ws.Range("D19") '13.16'
$ws.Range("E19").Value = '  +0.32%  '
Set-TextValue $ws.Range("D20") '6.62'
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("D21").Value = '0.0₃0943'
$ws.Range("E21").Value = '  +0.86%  '
Set-TextValue $ws.Range("D22") '2.72'
$ws.Range("E22").Value = '  +2.97%  '
Set-TextValue $ws.Range("D23") '72.21'
$ws.Range("E23").Value = '  +2.70%  '
Set-TextValue $ws.Range("D24") '275.59'
$ws.Range("E24").Value = '  +12.53%  '
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("E26").Value = '  -0.02%  '
Set-TextValue $ws.Range("D27") '25.90'
$ws.Range("E27").Value = '  +0.86%  '
Set-TextValue $ws.Range("D28") '2.39'
$ws.Range("E28").Value = '  +4.87%  '
Set-TextValue $ws.Range("D29") '10.05'
$ws.Range("E29").Value = '  +0.56%  '
Set-TextValue $ws.Range("D30") '0.140'
$ws.Range("E30").Value = '  +2.17%  '
Set-TextValue $ws.Range("D31") '35.39'
$ws.Range("E31").Value = '  +2.03%  '
Set-TextValue $ws.Range("D32") '49.53'
$ws.Range("E32").Value = '  -0.56%  '
Set-TextValue $ws.Range("D33") '19.35'
$ws.Range("E33").Value = '  -5.82%  '
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("E35").Value = '  -0.06%  '
Set-TextValue $ws.Range("D36") '0.0785'
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("E37").Value = '  -0.29%  '
$ws.Range("E38").Value = '  -2.06%  '
$ws.Range("E42").Value = '  -0.87%  '
Set-TextValue $ws.Range("D43") '21.65'
$ws.Range("E43").Value = '  -5.68%  '
$ws.Range("E44").Value = '  +3.23%  '
$ws.Range("D45").Value = '2.019.66'
$ws.Range("E45").Value = '  +1.14%  '
$ws.Range("E46").Value = '  +2.82%  '
$ws.Range("E47").Value = '  +4.64%  '
$ws.Range("E48").Value = '  -0.41%  '
Set-TextValue $ws.Range("D49") '9.03'
$ws.Range("E49").Value = '  -1.71%  '
Set-TextValue $ws.Range("D50") '5.19'
$ws.Range("E50").Value = '  +2.03%  '
Set-TextValue $ws.Range("D51") '79.56'
$ws.Range("E51").Value = '  +2.64%  '
